$d = $word.ActiveDocument

# --- 1. Title paragraph: "Regras de Comunicação" -> "REGRAS DE COMUNICAÇÃO" ---
#      also becomes bold, and shrinks from sz 32 (16pt) to sz 28 (14pt).
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Text = "REGRAS DE COMUNICAÇÃO"
$titlePara.Range.Font.Bold = 1
$titlePara.Range.Font.BoldBi = 1
$titlePara.Range.Font.Size = 14
$titlePara.Range.Font.SizeBi = 14

# --- 2. Merge the "Deixar o repositório do GitHub atualizado" runs ---
#      (collapses the GitHub proofErr spell-check split back into one run)
#      without swallowing the trailing " localmente;" run into the merge:
#      temporarily mark it bold so the run-coalescing pass after the
#      replace leaves it as its own <w:r>, then un-bold it again (plain
#      formatting edits do not trigger run coalescing).
$ghPara = $d.Paragraphs(2)
$ghFull = $ghPara.Range.Text
$ghStart = $ghPara.Range.Start
$tailText = " localmente;"
$tailIdx = $ghFull.IndexOf($tailText)
$tailRng = $d.Range($ghStart + $tailIdx, $ghStart + $tailIdx + $tailText.Length)
$tailRng.Font.Bold = 1

$d.Content.Find.Execute("Deixar o repositório do GitHub atualizado", $false, $false, $false, $false, $false, $true, 1, $false, "Deixar o repositório do GitHub atualizado", 2) | Out-Null

$ghPara2 = $d.Paragraphs(2)
$ghFull2 = $ghPara2.Range.Text
$ghStart2 = $ghPara2.Range.Start
$tailIdx2 = $ghFull2.IndexOf($tailText)
$tailRng2 = $d.Range($ghStart2 + $tailIdx2, $ghStart2 + $tailIdx2 + $tailText.Length)
$tailRng2.Font.Bold = 0

# --- 3. Split "Estabelecer horários..." so the _GoBack bookmark moves here ---
#      (placed right after "hor", before "ários e dias...") and drop it
#      from the end of the WhatsApp paragraph (a document only keeps one
#      _GoBack bookmark, so re-adding it moves it automatically).
$estPara = $d.Paragraphs(4)
$estText = $estPara.Range.Text
$splitAt = $estText.IndexOf("horá") + 4
$estStart = $estPara.Range.Start
$markRange = $d.Range($estStart + $splitAt, $estStart + $splitAt)
$d.Bookmarks.Add("_GoBack", $markRange) | Out-Null
